# Add a new time-log entry (row 16) to the hours tracking sheet and
# update the new "printing plots and adding titles" task description.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting from the row above (A15) so the new date cell
# (A16) keeps the same date number format (style) used by the rest of
# column A, then set its value.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A16").Value = 45736

# Hours spent and task description for the new entry.
$ws.Range("B16").Value = ".5 hours"
$ws.Range("C16").Value = "printing plots and adding titles"

# Match the saved selection state recorded in the workbook.
$ws.Range("C11").Select()
